# Apply the changes described by the diff:
#  - Add two new cell styles (border top+bottom, and border top+bottom+right)
#  - Apply those styles to the header spacer cells C1/D1 (and F1/G1 on sheet2)
#  - Rename "fedcore" -> "approach" in the header rows
#  - Remove the stray empty G5 cell on sheet2

$wb = $excel.ActiveWorkbook

# ---- Sheet "quality_comparison" ----
$ws1 = $wb.Worksheets.Item("quality_comparison")

$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$c1.Borders.Item(8).Weight = 2      # xlThin
$c1.Borders.Item(9).Weight = 2
$c1.Borders.Item(8).ColorIndex = -4105  # xlColorIndexAutomatic
$c1.Borders.Item(9).ColorIndex = -4105

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$d1.Borders.Item(8).Weight = 2
$d1.Borders.Item(10).Weight = 2
$d1.Borders.Item(9).Weight = 2
$d1.Borders.Item(8).ColorIndex = -4105
$d1.Borders.Item(10).ColorIndex = -4105
$d1.Borders.Item(9).ColorIndex = -4105

$ws1.Range("C2").Value = "approach"

# ---- Sheet "computational_comparison" ----
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1b.Borders.Item(8).LineStyle = 1
$c1b.Borders.Item(9).LineStyle = 1
$c1b.Borders.Item(8).Weight = 2
$c1b.Borders.Item(9).Weight = 2
$c1b.Borders.Item(8).ColorIndex = -4105
$c1b.Borders.Item(9).ColorIndex = -4105

$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1b.Borders.Item(8).LineStyle = 1
$d1b.Borders.Item(10).LineStyle = 1
$d1b.Borders.Item(9).LineStyle = 1
$d1b.Borders.Item(8).Weight = 2
$d1b.Borders.Item(10).Weight = 2
$d1b.Borders.Item(9).Weight = 2
$d1b.Borders.Item(8).ColorIndex = -4105
$d1b.Borders.Item(10).ColorIndex = -4105
$d1b.Borders.Item(9).ColorIndex = -4105

$f1b = $ws2.Range("F1")
$f1b.Style = "Normal"
$f1b.Borders.Item(8).LineStyle = 1
$f1b.Borders.Item(9).LineStyle = 1
$f1b.Borders.Item(8).Weight = 2
$f1b.Borders.Item(9).Weight = 2
$f1b.Borders.Item(8).ColorIndex = -4105
$f1b.Borders.Item(9).ColorIndex = -4105

$g1b = $ws2.Range("G1")
$g1b.Style = "Normal"
$g1b.Borders.Item(8).LineStyle = 1
$g1b.Borders.Item(10).LineStyle = 1
$g1b.Borders.Item(9).LineStyle = 1
$g1b.Borders.Item(8).Weight = 2
$g1b.Borders.Item(10).Weight = 2
$g1b.Borders.Item(9).Weight = 2
$g1b.Borders.Item(8).ColorIndex = -4105
$g1b.Borders.Item(10).ColorIndex = -4105
$g1b.Borders.Item(9).ColorIndex = -4105

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
